$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update progress value for row 15 ("1h") from 25 to 100
$ws.Range("D15").Value = 100
